$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2137518684603886
$ws.Range("C2").Value = 0.5171898355754858
$ws.Range("J2").Value = 0.06128550074738415
$ws.Range("O2").Value = 0.001494768310911809
$ws.Range("P2").Value = 0.1479820627802691
$ws.Range("S2").Value = 0.05829596412556054
$ws.Range("B3").Value = 0.01340482573726542
$ws.Range("C3").Value = 0.04825737265415549
$ws.Range("J3").Value = 0.1206434316353887
$ws.Range("P3").Value = 0.7372654155495979
$ws.Range("S3").Value = 0.08042895442359249
$ws.Range("J4").Value = 0.1348314606741573
$ws.Range("P4").Value = 0.7415730337078652
$ws.Range("S4").Value = 0.1235955056179775
$ws.Range("P5").Value = 0.875
$ws.Range("S5").Value = 0.125
$ws.Range("B6").Value = 0.06832298136645963
$ws.Range("D6").Value = 0.0124223602484472
$ws.Range("E6").Value = 0.004140786749482402
$ws.Range("F6").Value = 0.07453416149068323
$ws.Range("J6").Value = 0.3436853002070394
$ws.Range("O6").Value = 0.04968944099378882
$ws.Range("Q6").Value = 0.113871635610766
$ws.Range("R6").Value = 0.07660455486542443
$ws.Range("S6").Value = 0.2567287784679089
$ws.Range("B7").Value = 0.08823529411764706
$ws.Range("D7").Value = 0.0267379679144385
$ws.Range("E7").Value = 0.00267379679144385
$ws.Range("F7").Value = 0.06149732620320856
$ws.Range("J7").Value = 0.2192513368983957
$ws.Range("O7").Value = 0.0267379679144385
$ws.Range("Q7").Value = 0.1684491978609626
$ws.Range("R7").Value = 0.08021390374331551
$ws.Range("S7").Value = 0.3262032085561498
$ws.Range("B8").Value = 0.1054766734279919
$ws.Range("D8").Value = 0.01926977687626775
$ws.Range("E8").Value = 0.00101419878296146
$ws.Range("F8").Value = 0.05983772819472617
$ws.Range("J8").Value = 0.1825557809330629
$ws.Range("O8").Value = 0.02028397565922921
$ws.Range("Q8").Value = 0.1703853955375254
$ws.Range("R8").Value = 0.101419878296146
$ws.Range("S8").Value = 0.3397565922920893
$ws.Range("B9").Value = 0.1235294117647059
$ws.Range("D9").Value = 0.02058823529411765
$ws.Range("E9").Value = 0.002941176470588235
$ws.Range("F9").Value = 0.06176470588235294
$ws.Range("J9").Value = 0.1794117647058824
$ws.Range("O9").Value = 0.01176470588235294
$ws.Range("Q9").Value = 0.1588235294117647
$ws.Range("R9").Value = 0.1235294117647059
$ws.Range("S9").Value = 0.3176470588235294
$ws.Range("B10").Value = 0.09941720946177579
$ws.Range("D10").Value = 0.01816935207404868
$ws.Range("E10").Value = 0.001028453890983888
$ws.Range("F10").Value = 0.07062050051422694
$ws.Range("J10").Value = 0.2636270140555365
$ws.Range("O10").Value = 0.02639698320191978
$ws.Range("Q10").Value = 0.1895783339046966
$ws.Range("R10").Value = 0.07816249571477546
$ws.Range("S10").Value = 0.2529996571820363
$ws.Range("G11").Value = 0.1096196868008949
$ws.Range("J11").Value = 0.07829977628635347
$ws.Range("K11").Value = 0.1319910514541387
$ws.Range("L11").Value = 0.6644295302013423
$ws.Range("S11").Value = 0.01565995525727069
$ws.Range("G12").Value = 0.7884615384615384
$ws.Range("J12").Value = 0.1666666666666667
$ws.Range("L12").Value = 0.04166666666666666
$ws.Range("S12").Value = 0.003205128205128205
$ws.Range("F13").Value = 0.0119047619047619
$ws.Range("G13").Value = 0.6785714285714286
$ws.Range("J13").Value = 0.25
$ws.Range("S13").Value = 0.05952380952380952
$ws.Range("F15").Value = 0.02733485193621868
$ws.Range("H15").Value = 0.1708428246013667
$ws.Range("I15").Value = 0.05694760820045558
$ws.Range("J15").Value = 0.387243735763098
$ws.Range("K15").Value = 0.05466970387243736
$ws.Range("M15").Value = 0.01138952164009112
$ws.Range("N15").Value = 0.002277904328018223
$ws.Range("O15").Value = 0.04328018223234624
$ws.Range("S15").Value = 0.2460136674259681
$ws.Range("F16").Value = 0.02727272727272727
$ws.Range("H16").Value = 0.2022727272727273
$ws.Range("I16").Value = 0.07045454545454545
$ws.Range("J16").Value = 0.4022727272727273
$ws.Range("K16").Value = 0.1045454545454545
$ws.Range("M16").Value = 0.02045454545454545
$ws.Range("N16").Value = 0.002272727272727273
$ws.Range("O16").Value = 0.06136363636363636
$ws.Range("S16").Value = 0.1090909090909091
$ws.Range("F17").Value = 0.02386363636363636
$ws.Range("H17").Value = 0.1852272727272727
$ws.Range("I17").Value = 0.07045454545454545
$ws.Range("J17").Value = 0.4647727272727273
$ws.Range("K17").Value = 0.09545454545454546
$ws.Range("M17").Value = 0.01931818181818182
$ws.Range("N17").Value = 0.002272727272727273
$ws.Range("O17").Value = 0.07613636363636364
$ws.Range("S17").Value = 0.0625
$ws.Range("F18").Value = 0.01834862385321101
$ws.Range("H18").Value = 0.1995412844036697
$ws.Range("I18").Value = 0.09862385321100918
$ws.Range("J18").Value = 0.4793577981651376
$ws.Range("K18").Value = 0.06880733944954129
$ws.Range("M18").Value = 0.02064220183486239
$ws.Range("O18").Value = 0.04128440366972477
$ws.Range("S18").Value = 0.07339449541284404
$ws.Range("F19").Value = 0.02282157676348548
$ws.Range("H19").Value = 0.2282157676348548
$ws.Range("I19").Value = 0.07313278008298756
$ws.Range("J19").Value = 0.3973029045643153
$ws.Range("K19").Value = 0.0954356846473029
$ws.Range("M19").Value = 0.02074688796680498
$ws.Range("N19").Value = 0.0005186721991701245
$ws.Range("O19").Value = 0.06068464730290456
$ws.Range("S19").Value = 0.1011410788381743
